$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
$cl = $sm.CustomLayouts.Item(1)
$tcs = $cl.ThemeColorScheme
Write-Output "Count: $($tcs.Count)"
Write-Output $tcs.Colors(3).RGB
